# Normalize the vaccine/manufacturer labels in column A from all-caps
# (e.g. PFIZER_JANSSEN, JANSSEN, PFIZER_MODERNA, MODERNA, PFIZER) to
# title case (Pfizer_Janssen, Janssen, Pfizer_Moderna, Moderna, Pfizer).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$map = @{
    "PFIZER_JANSSEN" = "Pfizer_Janssen"
    "JANSSEN"        = "Janssen"
    "PFIZER_MODERNA" = "Pfizer_Moderna"
    "MODERNA"        = "Moderna"
    "PFIZER"         = "Pfizer"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
